# Apply targeted text fixes to the "card copy" document.
# Each replacement below corresponds to one <w:t> run change in the diff.
# We use literal (non-wildcard) Find/Replace since every search string is
# unique within the document, so a single Execute() call is sufficient and
# safe (no accidental multi-match edits).

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $new, 2, $false, $false, $false, $false) | Out-Null
}

Replace-Text "Початок будівництва: ????" "Початок будівництва: 2005"
Replace-Text "Завершення будівництва: ????" "Завершення будівництва: 2009"
Replace-Text "На початок будівництва: ????" "На початок будівництва: ____"
Replace-Text "На жовтень 2021р: ????" "На жовтень 2021р: ____"
Replace-Text "Комірок: ????" "Комірок: ____"
Replace-Text "Не житлових приміщень: ????" "Не житлових приміщень: ____"
Replace-Text "Середня вартість: ????" "Середня вартість: ____"
Replace-Text "Загальна площа квартир: ???" "Загальна площа квартир: 0.0"
Replace-Text "Площа вбудовано- прибудованих: ???" "Площа вбудовано- прибудованих: 0.0"
Replace-Text "Поверховість: 421" "Поверховість: 0-0"
Replace-Text "Комерційні приміщення: ???" "Комерційні приміщення: 1.0"
Replace-Text "Паркінг / к-ть місць: ???" "Паркінг / к-ть місць: 213"
Replace-Text "Гостьовий паркінг/ к-ть місць: ???" "Гостьовий паркінг/ к-ть місць: None"
Replace-Text "Кількість …-кім.: 6" "Кількість 4-кім. і більше: 6"
Replace-Text "Площа 1-кім.: ???" "Площа 1-кім.: 0.0-1.0"
Replace-Text "Площа 3-кім.: ???" "Площа 3-кім.: 0.0-2.0"
Replace-Text "Площа 2-кім.: ???" "Площа 2-кім.: 0.0-2.0"
Replace-Text "Площа …-кім.: ???" "Площа 4-кім. і більше: 0.0-4.0"

Write-Output "Done applying replacements."
